$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2-270). The update bumps that date by one day (45178 -> 45179)
# across the whole column.
$range = $ws.Range("C2:C270")
$range.Value = 45179
